$d = $word.ActiveDocument

# --- Step 1: paragraph 15 ("Les contrôles doivent être facile...") ---
# Only the paragraph MARK formatting (size) changes from 11pt (22 half-pts) to 12pt (24 half-pts).
# The paragraph's own run text/content is untouched.
$p15 = $d.Paragraphs(15)
$p15.Range.ParagraphFormat.Reset | Out-Null
$markRange = $d.Range($p15.Range.End - 1, $p15.Range.End)
$markRange.Font.Size = 12

# --- Step 2: paragraph 16 ("Sur Game Over, on doit pouvoir recommencer facilement le jeu.") ---
# Collapse its three runs into a single run with new text, 12pt, not italic.
$rsquo = [char]0x2019
$p16 = $d.Paragraphs(16)
$body16 = $d.Range($p16.Range.Start, $p16.Range.End - 1)
$body16.Text = "L" + $rsquo + "objectif du jeu est clair.  (10%)"
$body16 = $d.Range($p16.Range.Start, $p16.Range.End - 1)
$body16.Font.Size = 12

# --- Step 3: paragraph 17 (old "L'objectif du jeu est clair.  (10%)") ---
# Keep its paragraph/run formatting (12pt), only change the text.
$p17 = $d.Paragraphs(17)
$body17 = $d.Range($p17.Range.Start, $p17.Range.End - 1)
$body17.Text = "Utilisation judicieuse des menus.  (15%)"

# --- Steps 4-6: remove the three obsolete bullet paragraphs that followed ---
# After step 3 these are still paragraphs 18, 19 and 20:
#   18: "Utilisation des bonnes pratiques de jeu en VR, incluant les contrôles et interactions.  (25%)"
#   19: "Utilisation judicieuse des menus.  (15%)" (now a duplicate, to be removed)
#   20: "Recommencement facile sur Game Over.  (10%) "
# Delete paragraphs 18..20 inclusive (including their paragraph marks) in one shot.
$delStart = $d.Paragraphs(18).Range.Start
$delEnd = $d.Paragraphs(20).Range.End
$d.Range($delStart, $delEnd).Delete()
